$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 21 de Agosto de 2020 a las 23:52'
$ws.Range("B4").Value = 5789672
$ws.Range("C4").Value = 43400
$ws.Range("D4").Value = 3114822
$ws.Range("E4").Value = 2495835
$ws.Range("G4").Value = 959
$ws.Range("H4").Value = 179015
$ws.Range("B5").Value = 3532330
$ws.Range("C5").Value = 27233
$ws.Range("E5").Value = 765565
$ws.Range("G5").Value = 935
$ws.Range("H5").Value = 113358
$ws.Range("B8").Value = 603338
$ws.Range("C8").Value = 3398
$ws.Range("D8").Value = 500102
$ws.Range("E8").Value = 90393
$ws.Range("G8").Value = 225
$ws.Range("H8").Value = 12843
$ws.Range("B54").Value = 48661
$ws.Range("C54").Value = 358
$ws.Range("D54").Value = 45166
$ws.Range("E54").Value = 3314
$ws.Range("E76").Value = 7089
$ws.Range("G76").Value = 2
$ws.Range("H76").Value = 122
$ws.Range("B78").Value = 17310
$ws.Range("C78").Value = 61
$ws.Range("D78").Value = 14820
$ws.Range("E78").Value = 2378
$ws.Range("B81").Value = 15131
$ws.Range("C81").Value = 169
$ws.Range("D81").Value = 10282
$ws.Range("E81").Value = 4310
$ws.Range("G81").Value = 7
$ws.Range("H81").Value = 539
$ws.Range("B93").Value = 8932
$ws.Range("C93").Value = 56
$ws.Range("D93").Value = 7648
$ws.Range("E93").Value = 1231
$ws.Range("B99").Value = 8016
$ws.Range("C99").Value = 19
$ws.Range("E99").Value = 2373
$ws.Range("B103").Value = 6885
$ws.Range("C103").Value = 37
$ws.Range("D103").Value = 6186
$ws.Range("E103").Value = 541
$ws.Range("B107").Value = 5322
$ws.Range("C107").Value = 40
$ws.Range("D107").Value = 2929
$ws.Range("E107").Value = 2227
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 166
$ws.Range("A118").Value = 'Cabo Verde'
$ws.Range("B118").Value = 3412
$ws.Range("C118").Value = 44
$ws.Range("D118").Value = 2498
$ws.Range("E118").Value = 877
$ws.Range("H118").Value = 37
$ws.Range("A119").Value = 'Tailandia'
$ws.Range("B119").Value = 3390
$ws.Range("C119").Value = 1
$ws.Range("D119").Value = 3219
$ws.Range("E119").Value = 113
$ws.Range("H119").Value = 58
$ws.Range("A122").Value = 'Mayotte'
$ws.Range("B122").Value = 3237
$ws.Range("C122").Value = 16
$ws.Range("D122").Value = 2964
$ws.Range("E122").Value = 234
$ws.Range("H122").Value = 39
$ws.Range("A123").Value = 'Eslovaquia'
$ws.Range("B123").Value = 3225
$ws.Range("C123").Value = 123
$ws.Range("D123").Value = 2045
$ws.Range("E123").Value = 1147
$ws.Range("H123").Value = 33
$ws.Range("A124").Value = 'Mozambique'
$ws.Range("B124").Value = 3195
$ws.Range("C124").Value = 80
$ws.Range("D124").Value = 1406
$ws.Range("E124").Value = 1769
$ws.Range("H124").Value = 20
$ws.Range("B126").Value = 2780
$ws.Range("C126").Value = 63
$ws.Range("D126").Value = 1712
$ws.Range("E126").Value = 1057
$ws.Range("B127").Value = 2688
$ws.Range("C127").Value = 21
$ws.Range("D127").Value = 2007
$ws.Range("E127").Value = 556
$ws.Range("B132").Value = 2437
$ws.Range("C132").Value = 36
$ws.Range("D132").Value = 455
$ws.Range("E132").Value = 1898
$ws.Range("G132").Value = 3
$ws.Range("H132").Value = 84
$ws.Range("A136").Value = 'Siria'
$ws.Range("B136").Value = 2073
$ws.Range("C136").Value = 65
$ws.Range("D136").Value = 475
$ws.Range("E136").Value = 1515
$ws.Range("G136").Value = 1
$ws.Range("H136").Value = 83
$ws.Range("B137").Value = 2068
$ws.Range("C137").Value = 24
$ws.Range("D137").Value = 804
$ws.Range("E137").Value = 1170
$ws.Range("G137").Value = 1
$ws.Range("H137").Value = 94
$ws.Range("A138").Value = 'Islandia'
$ws.Range("B138").Value = 2050
$ws.Range("C138").Value = 10
$ws.Range("D138").Value = 1920
$ws.Range("E138").Value = 120
$ws.Range("H138").Value = 10
$ws.Range("B140").Value = 1906
$ws.Range("C140").Value = 7
$ws.Range("E140").Value = 306
$ws.Range("G140").Value = 1
$ws.Range("H140").Value = 542
$ws.Range("B154").Value = 1285
$ws.Range("C154").Value = 1
$ws.Range("E154").Value = 400
$ws.Range("B155").Value = 1239
$ws.Range("C155").Value = 27
$ws.Range("D155").Value = 891
$ws.Range("E155").Value = 321
$ws.Range("B156").Value = 1172
$ws.Range("C156").Value = 3
$ws.Range("D156").Value = 1083
$ws.Range("E156").Value = 20
$ws.Range("B162").Value = 891
$ws.Range("C162").Value = 3
$ws.Range("D162").Value = 830
$ws.Range("E162").Value = 46
